# Minnesota roster: fix row ordering that had gotten scrambled.
# Rows 4/5 swap (Jaylen Nowell <-> Naz Reid), rows 10/11 swap
# (Nathan Knight <-> Jordan McLaughlin), and rows 16/17/18 rotate
# (Josh Minott -> Mike Conley -> Nickeil Alexander-Walker -> Josh Minott).
# Only the cells whose displayed value actually changes are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Jaden McDaniels' position label (collateral fix alongside the reorder)
$ws.Range("D3").Value = "SF"

# Row 4: now Naz Reid (previously Jaylen Nowell)
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = "Naz Reid"
$ws.Range("D4").Value = "C"
$ws.Range("E4").Value = "6-9"
$ws.Range("F4").Value = 264
$ws.Range("G4").Value = "August 26, 1999"
$ws.Range("J4").Value = "LSU"
$ws.Range("K4").Value = "https://www.basketball-reference.com/players/r/reidna01.html"

# Row 5: now Jaylen Nowell (previously Naz Reid)
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "Jaylen Nowell"
$ws.Range("D5").Value = "SG"
$ws.Range("E5").Value = "6-4"
$ws.Range("F5").Value = 201
$ws.Range("G5").Value = "July 9, 1999"
$ws.Range("J5").Value = "Washington"
$ws.Range("K5").Value = "https://www.basketball-reference.com/players/n/nowelja01.html"

# Row 10: now Jordan McLaughlin (previously Nathan Knight)
$ws.Range("B10").Value = 6
$ws.Range("C10").Value = "Jordan McLaughlin"
$ws.Range("D10").Value = "PG"
$ws.Range("E10").Value = "5-11"
$ws.Range("F10").Value = 185
$ws.Range("G10").Value = "April 9, 1996"
$ws.Range("I10").Value = "'3"
$ws.Range("J10").Value = "USC"
$ws.Range("K10").Value = "https://www.basketball-reference.com/players/m/mclaujo01.html"

# Row 11: now Nathan Knight (previously Jordan McLaughlin)
$ws.Range("B11").Value = 13
$ws.Range("C11").Value = "Nathan Knight"
$ws.Range("D11").Value = "PF"
$ws.Range("E11").Value = "6-10"
$ws.Range("F11").Value = 253
$ws.Range("G11").Value = "September 20, 1997"
$ws.Range("I11").Value = "'2"
$ws.Range("J11").Value = "William & Mary"
$ws.Range("K11").Value = "https://www.basketball-reference.com/players/k/knighna01.html"

# Row 16: now Mike Conley (previously Josh Minott)
$ws.Range("B16").Value = 10
$ws.Range("C16").Value = "Mike Conley"
$ws.Range("D16").Value = "PG"
$ws.Range("E16").Value = "6-1"
$ws.Range("F16").Value = 175
$ws.Range("G16").Value = "October 11, 1987"
$ws.Range("I16").Value = "'15"
$ws.Range("J16").Value = "Ohio State"
$ws.Range("K16").Value = "https://www.basketball-reference.com/players/c/conlemi01.html"

# Row 17: now Nickeil Alexander-Walker (previously Mike Conley)
$ws.Range("B17").Value = 9
$ws.Range("C17").Value = "Nickeil Alexander-Walker"
$ws.Range("D17").Value = "SG"
$ws.Range("E17").Value = "6-6"
$ws.Range("F17").Value = 205
$ws.Range("G17").Value = "September 2, 1998"
$ws.Range("H17").Value = "ca"
$ws.Range("I17").Value = "'3"
$ws.Range("J17").Value = "Virginia Tech"
$ws.Range("K17").Value = "https://www.basketball-reference.com/players/a/alexani01.html"

# Row 18: now Josh Minott (previously Nickeil Alexander-Walker)
$ws.Range("B18").Value = 8
$ws.Range("C18").Value = "Josh Minott"
$ws.Range("D18").Value = "SF"
$ws.Range("E18").Value = "6-8"
$ws.Range("G18").Value = "November 25, 2002"
$ws.Range("H18").Value = "us"
$ws.Range("I18").Value = "R"
$ws.Range("J18").Value = "Memphis"
$ws.Range("K18").Value = "https://www.basketball-reference.com/players/m/minotjo01.html"
